$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 6256
$ws.Range("L3").Value = 6783
$ws.Range("J4").Value = 1883
$ws.Range("L4").Value = 1681
$ws.Range("L6").Value = 5572
$ws.Range("J7").Value = 29360
$ws.Range("L7").Value = 20692

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L2").Value = 70
$ws.Range("L7").Value = 233

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 486
$ws.Range("L4").Value = 93
$ws.Range("L7").Value = 1367

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 144
$ws.Range("L3").Value = 182
$ws.Range("L7").Value = 453

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L6").Value = 266
$ws.Range("L7").Value = 932

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 96
$ws.Range("L7").Value = 292

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L6").Value = 207
$ws.Range("L7").Value = 794

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L6").Value = 99
$ws.Range("L7").Value = 404

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L4").Value = 74
$ws.Range("L6").Value = 166
$ws.Range("L7").Value = 660
$ws.Range("L8").Value = 1367
$ws.Range("L9").Value = 120
$ws.Range("L11").Value = 341
$ws.Range("L14").Value = 100
$ws.Range("L16").Value = 48
$ws.Range("L19").Value = 559
$ws.Range("L23").Value = 219
$ws.Range("L27").Value = 178
$ws.Range("L29").Value = 1154
$ws.Range("L33").Value = 932
$ws.Range("L34").Value = 113
$ws.Range("L37").Value = 794
$ws.Range("L42").Value = 658
$ws.Range("L43").Value = 154
$ws.Range("L44").Value = 140
$ws.Range("L48").Value = 272
$ws.Range("L49").Value = 113
$ws.Range("L51").Value = 257
$ws.Range("L52").Value = 441
$ws.Range("L53").Value = 233
$ws.Range("L55").Value = 221
$ws.Range("L60").Value = 139
$ws.Range("J63").Value = 236
$ws.Range("L63").Value = 64
$ws.Range("L65").Value = 404
$ws.Range("L66").Value = 59
$ws.Range("L67").Value = 720
$ws.Range("L70").Value = 36
$ws.Range("L77").Value = 137
$ws.Range("L78").Value = 273
$ws.Range("L80").Value = 68
$ws.Range("L83").Value = 453
$ws.Range("L84").Value = 199
$ws.Range("L85").Value = 1028
$ws.Range("L90").Value = 218
$ws.Range("L91").Value = 279
$ws.Range("L92").Value = 63
$ws.Range("L94").Value = 253
$ws.Range("L95").Value = 292
$ws.Range("J101").Value = 29360
$ws.Range("L101").Value = 20692

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 204
$ws.Range("L6").Value = 167
$ws.Range("L7").Value = 720

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L2").Value = 67
$ws.Range("L7").Value = 199

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 347
$ws.Range("L4").Value = 63
$ws.Range("L7").Value = 1154

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 111
$ws.Range("L7").Value = 272

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 203
$ws.Range("L7").Value = 559

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L3").Value = 40
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 100

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L3").Value = 50
$ws.Range("L7").Value = 166

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 179
$ws.Range("L7").Value = 658

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 71
$ws.Range("L3").Value = 90
$ws.Range("L4").Value = 32
$ws.Range("L7").Value = 273

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 221

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L6").Value = 55
$ws.Range("L7").Value = 219

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 96
$ws.Range("L7").Value = 279

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L6").Value = 158
$ws.Range("L7").Value = 660

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L6").Value = 94
$ws.Range("L7").Value = 253

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 97
$ws.Range("L6").Value = 87
$ws.Range("L7").Value = 341

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 120

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L2").Value = 28
$ws.Range("L7").Value = 63

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L2").Value = 49
$ws.Range("L7").Value = 178

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 73
$ws.Range("L7").Value = 218

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 80
$ws.Range("L7").Value = 257

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 139

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L6").Value = 47
$ws.Range("L7").Value = 154

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 428
$ws.Range("L6").Value = 213
$ws.Range("L7").Value = 1028

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 137

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L4").Value = 9
$ws.Range("L7").Value = 68

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 138
$ws.Range("L3").Value = 138
$ws.Range("L7").Value = 441

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L4").Value = 5
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L2").Value = 8
$ws.Range("L7").Value = 48
